$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price strings that look like plain decimal numbers (e.g. "28.94") would be
# auto-converted to numeric values by Excel when assigned via .Value. Prefixing
# with a leading apostrophe keeps them as text, exactly like typing '28.94 into
# the cell by hand, without touching the cell's number format.

$ws.Range('D2').Value = '29.849.15'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '1.639.80'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '''215.38'
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '''28.94'
$ws.Range('E8').Value = '  -2.49%  '
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('E11').Value = '  -1.68%  '
$ws.Range('D12').Value = '1.874.52'
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('D13').Value = '1.639.49'
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('D14').Value = '''0.593'
$ws.Range('E14').Value = '  +3.87%  '
$ws.Range('D15').Value = '''9.55'
$ws.Range('E15').Value = '  +7.67%  '
$ws.Range('D16').Value = '''3.91'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').Value = '29.849.07'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').Value = '''64.36'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').Value = '''237.82'
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').Value = '''0.999'
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('D22').Value = '''9.93'
$ws.Range('E22').Value = '  +3.31%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  +2.33%  '
$ws.Range('D25').Value = '''157.55'
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').Value = '''15.60'
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('D28').Value = '''6.67'
$ws.Range('E28').Value = '  +0.96%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('E30').Value = '  +1.32%  '
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('D34').Value = '1.419.94'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('E35').Value = '  +2.43%  '
$ws.Range('E36').Value = '  -1.12%  '
$ws.Range('E37').Value = '  +1.83%  '
$ws.Range('E38').Value = '  -7.23%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '''0.571'
$ws.Range('E40').Value = '  +2.36%  '
$ws.Range('D41').Value = '''76.50'
$ws.Range('E41').Value = '  +10.53%  '
$ws.Range('D42').Value = '''0.0504'
$ws.Range('E42').Value = '  +0.38%  '
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('D44').Value = '''1.93'
$ws.Range('E44').Value = '  -2.71%  '
$ws.Range('E46').Value = '  -2.20%  '
$ws.Range('D47').Value = '''50.46'
$ws.Range('E47').Value = '  -7.65%  '
$ws.Range('D48').Value = '1.783.24'
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('E49').Value = '  -1.35%  '
$ws.Range('D50').Value = '''93.83'
$ws.Range('E50').Value = '  +5.79%  '
$ws.Range('E51').Value = '  +1.37%  '
